# EV usage/size probabilities implementation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sim_N: number of stochastic demand curves (row 7) ---
$ws.Range("D7").Value = 100

# --- flex_type description row (row 25) height changes automatically isn't guaranteed;
#     set explicitly to match target layout ---
$ws.Rows.Item(25).RowHeight = 195

# --- EV_size row (row 30): now holds size probabilities ---
$ws.Range("B30").Value = "Probabilités de taille du véhicule"
$ws.Range("D30").Value = "0.3, 0.5, 0.2"
$ws.Range("G30").Value = "[small, medium, large]"

# --- Insert a new row before the old EV_charger_power row so that EV_usage gets
#     its own row (row 31), pushing EV_charger_power / blank rows down ---
$ws.Rows.Item(31).Insert()

# --- EV_usage row (row 31) ---
$ws.Range("A31").Value = "EV_usage"
$ws.Range("B31").Value = "Probabilités de type d'usage du véhicule"
$ws.Range("D31").Value = "0.2, 0.5, 0.3"
$ws.Range("G31").Value = "[short, normal, long, int: (km/year)]"
$ws.Range("A31:G31").Style = $ws.Range("A30:G30").Style
$ws.Range("D31").Style = $ws.Range("D30").Style
$ws.Rows.Item(31).RowHeight = 19.5

# --- EV_km_per_year row (row 32, previously EV_usage) ---
$ws.Range("A32").Value = "EV_km_per_year"
$ws.Range("B32").Value = "Nombre de kilomètres par an (à la place de EV_usage)"
$ws.Range("D32").Value = 0
$ws.Range("G32").Value = "Si <=0, ne prend pas en compte cet input, simulation en fonciton de EV_usage"

# --- EV_charger_power row (row 33, previously blank spacer row) ---
$ws.Range("A33").Value = "EV_charger_power"
$ws.Range("B33").Value = "Puissance de charge du chargeur [kW]"
$ws.Range("D33").Value = 4
$ws.Range("A33:G33").Style = $ws.Range("A32:G32").Style
$ws.Range("D33").Style = $ws.Range("D31").Style
$ws.Rows.Item(33).RowHeight = 19.5

# --- Old Plot header (row 34) content cleared; row 35 (plt_plot) content cleared ---
$ws.Range("B34").Value = $null
$ws.Range("A35").Value = $null
$ws.Range("B35").Value = $null
$ws.Range("D35").Value = $null

# --- Insert two new blank rows before the Plot section so it moves from rows 34-35 to 37-38 ---
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).Insert()

$ws.Range("A36:G36").Style = $ws.Range("A33:G33").Style
$ws.Range("D36").Style = $ws.Range("D34").Style
$ws.Rows.Item(36).RowHeight = 18.75

$ws.Range("A37:G37").Style = $ws.Range("A34:G34").Style
$ws.Rows.Item(37).RowHeight = 19.5
$ws.Range("B37").Value = "Plot"

$ws.Range("A38:G38").Style = $ws.Range("A35:G35").Style
$ws.Rows.Item(38).RowHeight = 19.5
$ws.Range("A38").Value = "plt_plot"
$ws.Range("B38").Value = "Make a interactive plot"
$ws.Range("D38").Value = "True"

# --- Column widths (approximate autofit result) ---
$ws.Columns.Item(2).ColumnWidth = 44.46
$ws.Columns.Item(7).ColumnWidth = 63.03
